# feat: add 2022-Q3 data
#
# The workbook has three sheets: 总计 (summary), 2021-Q4, 2020-Q4.
# This adds a new "2022-Q3" sheet (with its own fund-holding data) right
# after "总计", pushing "2021-Q4" and "2020-Q4" down, and updates the
# "总计" summary sheet with a new row for the 2022-Q3 quarter.

$wb = $excel.ActiveWorkbook

# Helper: write a value into a range but force it to be stored as TEXT
# (mirrors the source file, where numeric-looking strings like "501307"
# or "0.15" are stored as strings, not numbers), without leaving a
# lingering explicit cell style behind.
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# ------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating the existing
#    "2021-Q4" sheet (so it inherits the same column layout / styling)
#    and placing the copy immediately before it.
# ------------------------------------------------------------------
$existing2021 = $wb.Worksheets.Item("2021-Q4")
$existing2021.Copy($existing2021, $null)

$q3 = $wb.Worksheets.Item("2021-Q4 (2)")
$q3.Name = "2022-Q3"

# Overwrite the copied data with the new 2022-Q3 fund-holding figures.
Set-TextValue $q3.Range("B2") "501307"
Set-TextValue $q3.Range("C2") "银河中证沪港深高股息指数（LOF）A"
Set-TextValue $q3.Range("D2") "0.15"
Set-TextValue $q3.Range("E2") "90.33"
Set-TextValue $q3.Range("F2") "1.27"
Set-TextValue $q3.Range("G2") "0.0019"
$q3.Range("H2").Value = 9

Set-TextValue $q3.Range("B3") "501308"
Set-TextValue $q3.Range("C3") "银河中证沪港深高股息指数（LOF）C"
Set-TextValue $q3.Range("D3") "0.01"
Set-TextValue $q3.Range("E3") "90.33"
Set-TextValue $q3.Range("F3") "1.27"
Set-TextValue $q3.Range("G3") "0.0001"
$q3.Range("H3").Value = 9

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q3
#    above the existing 2021-Q4 / 2020-Q4 rows.
# ------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Make room for a new row 4 by copying row 3's formatting into it,
# then fill rows 2-4 with their final values.
$zj.Range("A3").Copy($zj.Range("A4"))

$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2020-Q4"
$zj.Range("C4").Value = 2
$zj.Range("D4").Value = 0.03

$zj.Range("B3").Value = "2021-Q4"
$zj.Range("C3").Value = 2
$zj.Range("D3").Value = 0.03

$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 2
$zj.Range("D2").Value = 0

Write-Output "Done"
